$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the latest month row (row 22, 2025-09) with refreshed stats
$ws.Range("B22").Value = 6294
$ws.Range("C22").Value = 996
$ws.Range("D22").Value = 5840806
$ws.Range("E22").Value = 927.9958690816651
$ws.Range("F22").Value = 8.349113444654854
$ws.Range("G22").Value = 4.184100418410042
$ws.Range("H22").Value = 27.01932863020728
